# Add additional CI columns (LowerCISunXu, UpperCISunXu) to TableAUC table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths for newly introduced / resized columns ---
# (Target OOXML stored widths are 17 / 9.7265625 / 6.90625 / 10.90625 /
#  12.81640625 / 12.90625 respectively. This runtime's ColumnWidth setter
#  quantizes the stored width to the nearest 1/6th of a character, with a
#  fixed +5/6 offset between the "ColumnWidth" character value and the
#  stored width, so we request the character widths that land closest to
#  those targets after quantization.)
$ws.Columns.Item(3).ColumnWidth = 16.166666666666668
$ws.Columns.Item(4).ColumnWidth = 8.833333333333334
$ws.Columns.Item(5).ColumnWidth = 6.0
$ws.Columns.Item(6).ColumnWidth = 10.0
$ws.Columns.Item(7).ColumnWidth = 12.0
$ws.Columns.Item(8).ColumnWidth = 12.0

# --- Header row new column values (cols G and H) ---
$ws.Cells.Item(1, 7).Value2 = "LowerCISunXu"
$ws.Cells.Item(1, 8).Value2 = "UpperCISunXu"

# --- New CI data values for rows 2-46 (cols G and H) ---
$data = @(
    @(2, 0.96703708129600596, 0.98863636363636398),
    @(3, 0.94773891861905901, 0.98295454545454497),
    @(4, 0.93126722502248505, 0.97727272727272696),
    @(5, 0.92617971798033905, 0.97159090909090895),
    @(6, 0.91555891352403196, 0.96875),
    @(7, 0.91658326294811598, 0.96590909090909105),
    @(8, 0.91470351430631403, 0.96590909090909105),
    @(9, 0.91922821903416996, 0.96590909090909105),
    @(10, 0.91333102494881602, 0.96306818181818199),
    @(11, 0.90517577996444798, 0.96022727272727304),
    @(12, 0.90249063422686004, 0.95738636363636398),
    @(13, 0.88090215038608899, 0.94602272727272696),
    @(14, 0.87286804791819195, 0.94318181818181801),
    @(15, 0.87212016694806305, 0.94034090909090895),
    @(16, 0.86841477141568102, 0.9375),
    @(17, 0.84744358064512804, 0.92897727272727304),
    @(18, 0.84910535419679301, 0.92613636363636398),
    @(19, 0.831021247038912, 0.92045454545454497),
    @(20, 0.80994590681570899, 0.90625),
    @(21, 0.78122322107984998, 0.89488636363636398),
    @(22, 0.78356604995389401, 0.88920454545454497),
    @(23, 0.78409097024402497, 0.88636363636363602),
    @(24, 0.779296822968225, 0.88352272727272696),
    @(25, 0.75601360219948199, 0.87215909090909105),
    @(26, 0.73670988507617696, 0.86079545454545503),
    @(27, 0.73395860258469803, 0.85227272727272696),
    @(28, 0.72915000819590303, 0.84943181818181801),
    @(29, 0.72240747039615305, 0.84375),
    @(30, 0.71383991960457505, 0.83806818181818199),
    @(31, 0.71225281688962605, 0.83806818181818199),
    @(32, 0.70428519974447701, 0.83522727272727304),
    @(33, 0.69741811748195803, 0.82954545454545503),
    @(34, 0.68613744132913002, 0.82102272727272696),
    @(35, 0.65605827995773103, 0.79829545454545503),
    @(36, 0.64216787003992204, 0.78693181818181801),
    @(37, 0.62386647586039001, 0.77272727272727304),
    @(38, 0.60187985947063705, 0.76988636363636398),
    @(39, 0.56380644156330095, 0.73011363636363602),
    @(40, 0.54623943578998801, 0.71590909090909105),
    @(41, 0.48193798095135498, 0.66193181818181801),
    @(42, 0.48566764254681699, 0.66193181818181801),
    @(43, 0.42079728301720898, 0.61079545454545503),
    @(44, 0.4213244335886, 0.60511363636363602),
    @(45, 0.37610017877262097, 0.57954545454545503),
    @(46, 0.32154115660655902, 0.51136363636363602)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $g = $entry[1]
    $h = $entry[2]
    $ws.Cells.Item($row, 7).Value2 = $g
    $ws.Cells.Item($row, 8).Value2 = $h
}

# Match the number format (0.000) already used by columns D/E/F for the new G/H columns
$ws.Range("G2:H46").NumberFormat = "0.000"

# --- Header row formatting: center + wrap text ---
# (Apply cleanly to a single cell first via a named style, then propagate via
#  copy / paste-format so we don't leave behind intermediate orphaned styles.)
$headerStyle = $wb.Styles.Add("TableAUC_HeaderStyle")
$headerStyle.HorizontalAlignment = -4108
$headerStyle.WrapText = $true
$ws.Cells.Item(1, 1).Style = "TableAUC_HeaderStyle"
$ws.Cells.Item(1, 1).Copy()
$ws.Range("A1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Rows.Item(1).RowHeight = 43.5

# --- Selection update ---
$ws.Range("D8").Select()
